# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price/profit data updates to the Moogle_Profits workbook.
# Each sheet (crafting class) has per-leve rows with columns H..N holding market-price
# and profit figures pulled from an external data source; this script overwrites the
# specific cells whose source data changed, including a few cells that are newly
# populated (previously blank) and one cell that is cleared back to blank.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value2 = 6672.875
$ws.Range("I43").Value2 = 6812.35
$ws.Range("K43").Value2 = 6812.35
$ws.Range("M43").Value2 = -6743.35
$ws.Range("H53").Value2 = 1155.4
$ws.Range("I53").Value2 = 621
$ws.Range("J53").Value2 = 1689.8
$ws.Range("K53").Value2 = 621
$ws.Range("L53").Value2 = 1689.8
$ws.Range("M53").Value2 = 16
$ws.Range("N53").Value2 = -2963.8
$ws.Range("H92").Value2 = 788.7778
$ws.Range("I92").Value2 = 942.7143
$ws.Range("K92").Value2 = 942.7143
$ws.Range("M92").Value2 = 305.2857
$ws.Range("H100").Value2 = 1814.7778
$ws.Range("I100").Value2 = 1868.5883
$ws.Range("K100").Value2 = 1868.5883
$ws.Range("M100").Value2 = -1327.5883
$ws.Range("H107").Value2 = 491.91306
$ws.Range("I107").Value2 = 473.52942
$ws.Range("J107").Value2 = 544
$ws.Range("K107").Value2 = 473.52942
$ws.Range("L107").Value2 = 544
$ws.Range("M107").Value2 = 1446.47058
$ws.Range("N107").Value2 = -4384
$ws.Range("H116").Value2 = 7035.6113
$ws.Range("I116").Value2 = 6360.091
$ws.Range("J116").Value2 = 8097.143
$ws.Range("K116").Value2 = 6360.091
$ws.Range("L116").Value2 = 8097.143
$ws.Range("M116").Value2 = -2918.091
$ws.Range("N116").Value2 = -14981.143
$ws.Range("H126").Value2 = 80000
$ws.Range("J126").Value2 = 80000
$ws.Range("L126").Value2 = 80000
$ws.Range("N126").Value2 = -89880
$ws.Range("H132").Value2 = 3308.182
$ws.Range("I132").Value2 = 3255.3125
$ws.Range("K132").Value2 = 9765.9375
$ws.Range("M132").Value2 = -7235.9375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value2 = 55000
$ws.Range("J106").Value2 = 55000
$ws.Range("L106").Value2 = 55000
$ws.Range("N106").Value2 = -57524

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value2 = 1600.28
$ws.Range("I20").Value2 = 1745.238
$ws.Range("J20").Value2 = 839.25
$ws.Range("K20").Value2 = 1745.238
$ws.Range("L20").Value2 = 839.25
$ws.Range("M20").Value2 = -1498.238
$ws.Range("N20").Value2 = -1333.25
$ws.Range("H21").Value2 = 31632.334
$ws.Range("J21").Value2 = 31632.334
$ws.Range("L21").Value2 = 31632.334
$ws.Range("N21").Value2 = -32104.334
$ws.Range("H61").Value2 = 26590
$ws.Range("J61").Value2 = 26590
$ws.Range("L61").Value2 = 26590
$ws.Range("N61").Value2 = -27216
$ws.Range("H99").Value2 = 4819.091
$ws.Range("I99").Value2 = 2999.875
$ws.Range("J99").Value2 = 9670.333000000001
$ws.Range("K99").Value2 = 2999.875
$ws.Range("L99").Value2 = 9670.333000000001
$ws.Range("M99").Value2 = -1501.875
$ws.Range("N99").Value2 = -12666.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value2 = 19419.857
$ws.Range("I39").Value2 = 16506.5
$ws.Range("J39").Value2 = 36900
$ws.Range("K39").Value2 = 16506.5
$ws.Range("L39").Value2 = 36900
$ws.Range("M39").Value2 = -16115.5
$ws.Range("N39").Value2 = -37682
$ws.Range("H47").Value2 = 25000
$ws.Range("I47").Value2 = 20000
$ws.Range("J47").Value2 = 30000
$ws.Range("K47").Value2 = 20000
$ws.Range("L47").Value2 = 30000
$ws.Range("M47").Value2 = -19434
$ws.Range("N47").Value2 = -31132
$ws.Range("H49").Value2 = 19419.857
$ws.Range("I49").Value2 = 16506.5
$ws.Range("J49").Value2 = 36900
$ws.Range("K49").Value2 = 16506.5
$ws.Range("L49").Value2 = 36900
$ws.Range("M49").Value2 = -16324.5
$ws.Range("N49").Value2 = -37264
$ws.Range("H99").Value2 = 2065.0625
$ws.Range("I99").Value2 = 1900.9429
$ws.Range("J99").Value2 = 2506.923
$ws.Range("K99").Value2 = 1900.9429
$ws.Range("L99").Value2 = 2506.923
$ws.Range("M99").Value2 = -402.9429
$ws.Range("N99").Value2 = -5502.923
$ws.Range("H105").Value2 = 3470.1428
$ws.Range("I105").Value2 = 3258.4
$ws.Range("J105").Value2 = 3999.5
$ws.Range("K105").Value2 = 3258.4
$ws.Range("L105").Value2 = 3999.5
$ws.Range("M105").Value2 = -1511.4
$ws.Range("N105").Value2 = -7493.5
$ws.Range("H126").Value2 = 2065.0625
$ws.Range("I126").Value2 = 1900.9429
$ws.Range("J126").Value2 = 2506.923
$ws.Range("K126").Value2 = 5702.8287
$ws.Range("L126").Value2 = 7520.768999999999
$ws.Range("M126").Value2 = -3232.8287
$ws.Range("N126").Value2 = -12460.769
$ws.Range("H132").Value2 = 3906.36
$ws.Range("I132").Value2 = 2840.9524
$ws.Range("K132").Value2 = 8522.8572
$ws.Range("M132").Value2 = -5992.8572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 706.8570999999999
$ws.Range("I5").Value2 = 599.5
$ws.Range("J5").Value2 = 749.8
$ws.Range("K5").Value2 = 1798.5
$ws.Range("L5").Value2 = 2249.4
$ws.Range("M5").Value2 = -1686.5
$ws.Range("N5").Value2 = -2473.4
$ws.Range("H8").Value2 = 478.6154
$ws.Range("I8").Value2 = 478.6154
$ws.Range("K8").Value2 = 1435.8462
$ws.Range("M8").Value2 = -1296.8462
$ws.Range("H12").Value2 = 360.7857
$ws.Range("I12").Value2 = 72.8
$ws.Range("J12").Value2 = 520.7778
$ws.Range("K12").Value2 = 218.4
$ws.Range("L12").Value2 = 1562.3334
$ws.Range("M12").Value2 = -45.39999999999998
$ws.Range("N12").Value2 = -1908.3334
$ws.Range("H76").Value2 = 12751.363
$ws.Range("I76").Value2 = 8460.833000000001
$ws.Range("J76").Value2 = 17900
$ws.Range("K76").Value2 = 25382.499
$ws.Range("L76").Value2 = 53700
$ws.Range("M76").Value2 = -24999.499
$ws.Range("N76").Value2 = -54466
$ws.Range("H79").Value2 = 12751.363
$ws.Range("I79").Value2 = 8460.833000000001
$ws.Range("J79").Value2 = 17900
$ws.Range("K79").Value2 = 25382.499
$ws.Range("L79").Value2 = 53700
$ws.Range("M79").Value2 = -24056.499
$ws.Range("N79").Value2 = -56352
$ws.Range("H80").Value2 = 4665
$ws.Range("I80").Value2 = 4000
$ws.Range("J80").Value2 = 4997.5
$ws.Range("K80").Value2 = 12000
$ws.Range("L80").Value2 = 14992.5
$ws.Range("M80").Value2 = -11064
$ws.Range("N80").Value2 = -16864.5
$ws.Range("H83").Value2 = 4665
$ws.Range("I83").Value2 = 4000
$ws.Range("J83").Value2 = 4997.5
$ws.Range("K83").Value2 = 36000
$ws.Range("L83").Value2 = 44977.5
$ws.Range("M83").Value2 = -31320
$ws.Range("N83").Value2 = -54337.5
$ws.Range("H92").Value2 = 291.66666
$ws.Range("J92").Value2 = 291.66666
$ws.Range("L92").Value2 = 874.9999799999999
$ws.Range("N92").Value2 = -3370.99998
$ws.Range("H113").Value2 = 1279.3077
$ws.Range("J113").Value2 = 1675.2222
$ws.Range("L113").Value2 = 5025.6666
$ws.Range("N113").Value2 = -9365.6666
$ws.Range("H135").Value2 = 706.8570999999999
$ws.Range("I135").Value2 = 599.5
$ws.Range("J135").Value2 = 749.8
$ws.Range("K135").Value2 = 5395.5
$ws.Range("L135").Value2 = 6748.2
$ws.Range("M135").Value2 = -2860.5
$ws.Range("N135").Value2 = -11818.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value2 = 8246.157999999999
$ws.Range("I132").Value2 = 7804.5884
$ws.Range("J132").Value2 = 11999.5
$ws.Range("K132").Value2 = 23413.7652
$ws.Range("L132").Value2 = 35998.5
$ws.Range("M132").Value2 = -20883.7652
$ws.Range("N132").Value2 = -41058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value2 = 2459.6428
$ws.Range("I93").Value2 = 2461.2222
$ws.Range("J93").Value2 = 2456.8
$ws.Range("K93").Value2 = 2461.2222
$ws.Range("L93").Value2 = 2456.8
$ws.Range("M93").Value2 = -1213.2222
$ws.Range("N93").Value2 = -4952.8
$ws.Range("H105").Value2 = 16486.334
$ws.Range("J105").Value2 = 16450
$ws.Range("L105").Value2 = 16450
$ws.Range("N105").Value2 = -23438

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value2 = 419.8
$ws.Range("I100").Value2 = 442.57144
$ws.Range("K100").Value2 = 885.14288
$ws.Range("M100").Value2 = -344.14288
$ws.Range("H94").Value2 = 0
$ws.Range("J94").Value2 = 0
$ws.Range("L94").Value2 = 0
$ws.Range("N94").ClearContents()
